$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) transitLineToVehicle: insert 3 new rows (492-494) for the HSR "137_"
#    lines, just above the trailing repeated header row (which slides from
#    492 down to 495).
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("transitLineToVehicle")

$ws2.Range("A492:A494").EntireRow.Insert()

# Column A - line names
$ws2.Range("A492").Value = "137_A"
$ws2.Range("A493").Value = "137_B"
$ws2.Range("A494").Value = "137_C"

# Column B - System
$ws2.Range("B492").Value = "HSR"
$ws2.Range("B493").Value = "HSR"
$ws2.Range("B494").Value = "HSR"

# Column D - Line
$ws2.Range("D492").Value = "SFG"
$ws2.Range("D493").Value = "SFSJ"
$ws2.Range("D494").Value = "SJG"

# Column E - FullLineName
$ws2.Range("E492").Value = "SF - GILROY"
$ws2.Range("E493").Value = "SF - SAN JOSE"
$ws2.Range("E494").Value = "SAM JOSE - GILROY"

# Column F - AM VehicleType
$ws2.Range("F492").Value = "HSR6"
$ws2.Range("F493").Value = "HSR6"
$ws2.Range("F494").Value = "HSR4"

# Column G - PM VehicleType
$ws2.Range("G492").Value = "HSR6"
$ws2.Range("G493").Value = "HSR6"
$ws2.Range("G494").Value = "HSR4"

# Column H - OP Vehicle Type
$ws2.Range("H492").Value = "HSR10"
$ws2.Range("H493").Value = "HSR10"
$ws2.Range("H494").Value = "HSR10"

# Column C - Stripped (formula, same pattern used throughout the column) -
# set on the whole new block at once so it is stored as one shared formula,
# matching how the rest of the column is already laid out.
$ws2.Range("C492:C494").Formula = '=RIGHT($A492,LEN($A492)-FIND("_",$A492))'

# Visually mark row 492 as the first row of the new "HSR" system block, same
# way every other system boundary is highlighted in this sheet.
$ws2.Range("A32:H32").Copy()
$ws2.Range("A492:H492").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 2) transitPrefixToVehicle: the "137_" prefix now maps to vehicle type "HSR"
#    instead of "Unknown Train".
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("transitPrefixToVehicle")
$ws1.Range("C65").Value = "HSR"

# ---------------------------------------------------------------------------
# 3) transitVehicleToCapacity: add capacity rows for the three new HSR
#    vehicle types (4/6/10 train consists).
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("transitVehicleToCapacity")

$ws3.Range("A48").Value = "HSR4"
$ws3.Range("B48").Value = 1800
$ws3.Range("C48").Value = 1530
$ws3.Range("D48").Value = "HSR 4 trains"
$ws3.Range("E48").Value = 0
$ws3.Range("F48").Value = 0
$ws3.Range("G48").Value = 0
$ws3.Range("H48").Value = 0

$ws3.Range("A49").Value = "HSR6"
$ws3.Range("B49").Value = 2700
$ws3.Range("C49").Value = 2295
$ws3.Range("D49").Value = "HSR 6 trains"
$ws3.Range("E49").Value = 0
$ws3.Range("F49").Value = 0
$ws3.Range("G49").Value = 0
$ws3.Range("H49").Value = 0

$ws3.Range("A50").Value = "HSR10"
$ws3.Range("B50").Value = 4500
$ws3.Range("C50").Value = 3825
$ws3.Range("D50").Value = "HSR 10 trains"
$ws3.Range("E50").Value = 0
$ws3.Range("F50").Value = 0
$ws3.Range("G50").Value = 0
$ws3.Range("H50").Value = 0

# ---------------------------------------------------------------------------
# 4) Leave the cursor/selection where the editor ended up on each sheet.
#    transitVehicleToCapacity is selected last so it stays the active tab,
#    matching the workbook's original/active-tab state.
# ---------------------------------------------------------------------------
$ws1.Range("C13").Select()
$ws2.Range("H497").Select()
$ws3.Range("E53").Select()
